{"js": "// Insert four new paragraphs (\"Authour: Joe\", \"Location: Paris\",\n// \"Date June 9\", \"end\") at the very top of the document body, before the\n// existing \"This is a my first document.\" paragraph.\nconst body = context.document.body;\n\nlet p = body.insertParagraph(\"Authour: Joe\", Word.InsertLocation.start);\np = p.insertParagraph(\"Location: Paris\", Word.InsertLocation.after);\np = p.insertParagraph(\"Date June 9\", Word.InsertLocation.after);\np = p.insertParagraph(\"end\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert four new paragraphs (\"Authour: Joe\", \"Location: Paris\",\n# \"Date June 9\", \"end\") at the very top of the document, before the\n# existing \"This is a my first document.\" paragraph.\n$d = $word.ActiveDocument\n\n# Anchor on the range of what is currently the first paragraph, then\n# repeatedly insert text + a paragraph mark immediately before it. Each\n# InsertBefore call pushes the previous content further down, so the new\n# lines are added in reverse order to land in the desired final order.\n$r = $d.Paragraphs(1).Range\n$r.InsertBefore(\"end`r\")\n$r.InsertBefore(\"Date June 9`r\")\n$r.InsertBefore(\"Location: Paris`r\")\n$r.InsertBefore(\"Authour: Joe`r\")\n"}
